# Added single qubit classifier and analysis
# - Reorders the metric columns in the existing table (C:G) to
#   Mean val Acc, Best LR (Mean), Std Val loss, Max Val Acc, Best LR (Max val)
# - Adds a second table below (rows 10-17) titled "With crossvalidation"
#   holding the single-qubit classifier results with cross-validation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing table (rows 2-8): columns C-G re-ordered, values updated ---
$ws.Range("B2").Value = "Model"
$ws.Range("C2").Value = "Mean val Acc"
$ws.Range("D2").Value = "Best LR (Mean)"
$ws.Range("E2").Value = "Std Validation loss across all runs for best mean lr"
$ws.Range("F2").Value = "Max Val Acc"
$ws.Range("G2").Value = "Best LR (Max val)"

$ws.Range("B3").Value = "V"
$ws.Range("C3").Value = 0.46
$ws.Range("D3").Value = 0.05
$ws.Range("E3").Value = 0.7
$ws.Range("F3").Value = 0.72
$ws.Range("G3").Value = 0.1

$ws.Range("B4").Value = "W"
$ws.Range("C4").Value = 0.5
$ws.Range("D4").Value = 0.05
$ws.Range("E4").Value = 0.44
$ws.Range("F4").Value = 0.9
$ws.Range("G4").Value = 0.05

$ws.Range("B5").Value = "X"
$ws.Range("C5").Value = 0.52
$ws.Range("D5").Value = 0.01
$ws.Range("E5").Value = 0.12
$ws.Range("F5").Value = 0.84
$ws.Range("G5").Value = 0.01
$ws.Range("I5").Value = "Als de max en mean een andere beste LR hebben is het verschil tussen de means voor de twee lrs meestal klein (.01)"

$ws.Range("B6").Value = "Y"
$ws.Range("C6").Value = 0.45
$ws.Range("D6").Value = 0.05
$ws.Range("E6").Value = 0.21
$ws.Range("F6").Value = 0.78
$ws.Range("G6").Value = 0.05

$ws.Range("B7").Value = "Z"
$ws.Range("C7").Value = 0.44
$ws.Range("D7").Value = 0.01
$ws.Range("E7").Value = 0.21
$ws.Range("F7").Value = 0.89
$ws.Range("G7").Value = 0.05

$ws.Range("B8").Value = "E"
$ws.Range("C8").Value = 0.43
$ws.Range("D8").Value = 0.01
$ws.Range("E8").Value = 0.22
$ws.Range("F8").Value = 0.92
$ws.Range("G8").Value = 0.01

# --- New section: single qubit classifier with crossvalidation ---
$ws.Range("B11").Value = "Model"
$ws.Range("C11").Value = "Mean val Acc"
$ws.Range("D11").Value = "Best LR (Mean)"
$ws.Range("E11").Value = "Std Validation loss across all runs for best mean lr"
$ws.Range("F11").Value = "Max Val Acc"
$ws.Range("G11").Value = "Best LR (Max val)"

$ws.Range("B12").Value = "V"
$ws.Range("C12").Value = 0.72
$ws.Range("D12").Value = 0.05
$ws.Range("E12").Value = 0.45
$ws.Range("F12").Value = 0.95
$ws.Range("G12").Value = 0.05

$ws.Range("B13").Value = "W"
$ws.Range("C13").Value = 0.89
$ws.Range("D13").Value = 0.05
$ws.Range("E13").Value = 0.06
$ws.Range("F13").Value = 0.96

$ws.Range("B14").Value = "X"
$ws.Range("C14").Value = 0.88
$ws.Range("E14").Value = 0.06
$ws.Range("F14").Value = 0.98

$ws.Range("B15").Value = "Y"
$ws.Range("C15").Value = 0.53
$ws.Range("E15").Value = 0.08
$ws.Range("F15").Value = 0.86

$ws.Range("B16").Value = "Z"
$ws.Range("C16").Value = 0.88
$ws.Range("E16").Value = 0.07
$ws.Range("F16").Value = 0.96

$ws.Range("B17").Value = "A"
$ws.Range("C17").Value = 0.9
$ws.Range("E17").Value = 0.06
$ws.Range("F17").Value = 0.98

# Section title (set after the table body so shared-string insertion order
# mirrors the original authoring order)
$ws.Range("A10").Value = "With crossvalidation"

# Match the saved selection in the workbook (cell I5 was last selected)
$ws.Range("I5").Select()
